$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2,2,2,1,1,1,2,1,2,2,2,1,2,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2  # Column B is 2
    $ws.Cells.Item(2, $col).Value = $values[$i]
}
